# RequestCodes.xlsx update
# - Adds two new shared strings (used by the new row in sheet1)
# - Adds highlight formatting (copy of the style already used on A2/A8) to
#   A3:A5 in sheet1 and A7/A9/A10 in sheet2
# - Adds a new row (10) to sheet1 describing the new "Respuesta de union a
#   sala" -> (Selección | Solicitante / Boolean | string) request code
# - Updates the active selection on sheet1 and sets the page to portrait

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- sheet1: highlight column A on rows 3,4,5 (same look as A2/A8) ---
$ws1.Range("A2").Copy()
$ws1.Range("A3:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- sheet1: new row 10 ---
$ws1.Range("B10").Value = "Respuesta de union a sala"
$ws1.Range("C10").Value = 109
$ws1.Range("D10").Value = "Selección | Solicitante"
$ws1.Range("E10").Value = "Boolean | string"
$ws1.Range("F10").Value = "Individual (Requester)"

# --- sheet1: selection + page setup ---
$ws1.Range("E13").Select()
$ws1.PageSetup.Orientation = 1

# --- sheet2: highlight column A on rows 7,9,10 (same look as A2/A3/A6) ---
$ws2.Range("A2").Copy()
$ws2.Range("A7").PasteSpecial(-4122)
$ws2.Range("A9").PasteSpecial(-4122)
$ws2.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Sheet2 ("Client") was the originally active/visible tab - restore that
# after touching sheet1 so selecting E13 above doesn't steal the active tab.
$ws2.Activate()

$wb.Save()
